# ---------------------------------------------------------------------------
# Edit script: add "version" schema column + "description" column to the
# "Export as TSV" sheet, add corresponding comments and data validation,
# and insert the "version list" / "signal_type list" lookup sheets.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export as TSV")

$headers = @(
    'version',
    'description',
    'donor_id',
    'tissue_id',
    'execution_datetime',
    'protocols_io_doi',
    'operator',
    'operator_email',
    'pi',
    'pi_email',
    'assay_category',
    'assay_type',
    'analyte_class',
    'is_targeted',
    'acquisition_instrument_vendor',
    'acquisition_instrument_model',
    'preparation_instrument_vendor',
    'preparation_instrument_model',
    'section_prep_protocols_io_doi',
    'reagent_prep_protocols_io_doi',
    'number_of_channels',
    'number_of_sections',
    'ablation_distance_between_shots_x_value',
    'ablation_distance_between_shots_x_units',
    'ablation_distance_between_shots_y_value',
    'ablation_distance_between_shots_y_units',
    'ablation_frequency_value',
    'ablation_frequency_unit',
    'roi_description',
    'roi_id',
    'acquisition_id',
    'max_x_width_value',
    'max_x_width_unit',
    'max_y_height_value',
    'max_y_height_unit',
    'segment_data_format',
    'signal_type',
    'antibodies_path',
    'contributors_path',
    'data_path',
)

$comments = @(
    'Version of the schema to use when validating this metadata.',
    'Free-text description of this assay.',
    'HuBMAP Display ID of the donor of the assayed tissue.',
    'HuBMAP Display ID of the assayed tissue.',
    'Start date and time of assay, typically a date-time stamped folder generated by the acquisition instrument. YYYY-MM-DD hh:mm, where YYYY is the year, MM is the month with leading 0s, and DD is the day with leading 0s, hh is the hour with leading zeros, mm are the minutes with leading zeros.',
    'DOI for protocols.io referring to the protocol for this assay.',
    'Name of the person responsible for executing the assay.',
    'Email address for the operator.',
    'Name of the principal investigator responsible for the data.',
    'Email address for the principal investigator.',
    'Each assay is placed into one of the following 3 general categories: generation of images of microscopic entities, identification & quantitation of molecules by mass spectrometry, and determination of nucleotide sequence.',
    'The specific type of assay being executed.',
    'Analytes are the target molecules being measured with the assay.',
    'Specifies whether or not a specific molecule(s) is/are targeted for detection/measurement by the assay. The CODEX analyte is protein.',
    'An acquisition instrument is the device that contains the signal detection hardware and signal processing software. Assays generate signals such as light of various intensities or color or signals representing the molecular mass.',
    'Manufacturers of an acquisition instrument may offer various versions (models) of that instrument with different features or sensitivities. Differences in features or sensitivities may be relevant to processing or interpretation of the data.',
    'The manufacturer of the instrument used to prepare the sample for the assay.',
    'The model number/name of the instrument used to prepare the sample for the assay',
    'DOI for protocols.io referring to the protocol for preparing tissue sections for the assay.',
    'DOI for protocols.io referring to the protocol for preparing reagents for the assay.',
    'Number of mass channels measured',
    'Number of sections',
    'x resolution. Distance between laser ablation shots in the X-dimension.',
    'Units of x resolution distance between laser ablation shots.',
    'y resolution. Distance between laser ablation shots in the Y-dimension.',
    'Units of y resolution distance between laser ablation shots.',
    'Frequency value of laser ablation (in Hz)',
    'Frequency unit of laser ablation',
    'A description of the region of interest (ROI) captured in the image.',
    'Multiple images (1-n) are acquired from regions of interest (ROI1, ROI2, ROI3, etc) on a slide. The roi_id is a number from 1-n representing the ROI captured on a slide.',
    'The acquisition_id refers to the directory containing the ROI images for a slide. Together, the acquisition_id and the roi_id indicate the slide-ROI represented in the image.',
    'Image width value of the ROI acquisition',
    'Units of image width of the ROI acquisition',
    'Image height value of the ROI acquisition',
    'Units of image height of the ROI acquisition',
    'This refers to the data type, which is a "float" for the IMC counts.',
    'Type of signal measured per channel (usually dual counts)',
    'Relative path to file with antibody information for this dataset.',
    'Relative path to file with ORCID IDs for contributors for this dataset.',
    'Relative path to file or directory with instrument data. Downstream processing will depend on filename extension conventions.',
)

# ---------------------------------------------------------------------------
# 1. Remove the existing header-row comments (they are pinned to the old
#    column positions and do not move automatically when columns shift).
# ---------------------------------------------------------------------------
$oldLastCol = 38
for ($col = 1; $col -le $oldLastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    if ($cell.Comment -ne $null) {
        $cell.Comment.Delete()
    }
}

# ---------------------------------------------------------------------------
# 2. Insert two new columns at the front of the sheet: "version" and
#    "description". This shifts all existing data, formatting and data
#    validation ranges two columns to the right.
# ---------------------------------------------------------------------------
$ws.Range("A:B").Insert()

# ---------------------------------------------------------------------------
# 3. Rewrite the header row using the final column order, and re-apply the
#    header style (bold, centered, wrapped) to the two new cells.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Re-create all header comments (2 new + 38 shifted) at their final
#    column positions.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt $comments.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).AddComment($comments[$i])
}

# ---------------------------------------------------------------------------
# 5. Add the data validation list rule for the new "version" column.
# ---------------------------------------------------------------------------
$rngVersion = $ws.Range("A2:A1048576")
$rngVersion.Validation.Add(3, 1, 1, "'version list'!`$A`$1:`$A`$1")
$rngVersion.Validation.ErrorTitle = "Value must come from list"
$rngVersion.Validation.ErrorMessage = "Value must be one of: 1."
$rngVersion.Validation.ShowError = $true
$rngVersion.Validation.ShowInput = $true
$rngVersion.Validation.IgnoreBlank = $true

# ---------------------------------------------------------------------------
# 6. Insert the "version list" lookup sheet right before "assay_category
#    list", and the "signal_type list" lookup sheet at the very end.
# ---------------------------------------------------------------------------
$assayCategorySheet = $wb.Worksheets.Item("assay_category list")
$versionSheet = $wb.Worksheets.Add($assayCategorySheet)
$versionSheet.Name = "version list"
$versionSheet.Range("A1").Value = "1"

$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$signalSheet = $wb.Worksheets.Add($null, $lastSheet)
$signalSheet.Name = "signal_type list"
$signalSheet.Range("A1").Value = "dual count"
$signalSheet.Range("A2").Value = "pulse count"
$signalSheet.Range("A3").Value = "intensity value"

# Leave the original sheet selected/focused, matching the source file.
$ws.Activate()
$ws.Range("A1").Select()

